$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.127111911773682
$ws.Range("B1").Value = 1.608078002929688
$ws.Range("C1").Value = 4.589194297790527
$ws.Range("D1").Value = 0.5346905589103699
$ws.Range("E1").Value = 0.5934686660766602
